$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.083.16'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '3.421.64'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("D5").Value = "'579.45"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").Value = "'153.93"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = '  +5.13%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +1.63%  '
$ws.Range("E9").Value = '  +4.62%  '
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("E11").Value = '  +3.70%  '
$ws.Range("D12").Value = '4.008.16'
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").Value = "'28.65"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").Value = '3.425.25'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").Value = '62.109.22'
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("D18").Value = "'6.52"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = '  +2.92%  '
$ws.Range("D19").Value = "'14.42"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").Value = "'8.99"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("D21").Value = "'383.34"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = "'0.571"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = '  +2.33%  '
$ws.Range("D23").Value = "'75.94"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '3.562.94'
$ws.Range("E26").Value = '  -1.10%  '
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").Value = "'7.65"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = '  +1.53%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D31").Value = "'7.86"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = "'23.30"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("E35").Value = '  +4.89%  '
$ws.Range("E36").Value = '  +0.80%  '
$ws.Range("D37").Value = "'6.96"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("D38").Value = "'168.73"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").Value = "'30.83"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = '  -2.70%  '
$ws.Range("D40").Value = '3.458.49'
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("E41").Value = '  +1.85%  '
$ws.Range("E42").Value = '  +1.10%  '
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("E45").Value = '  -2.02%  '
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("D47").Value = '2.558.40'
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").Value = "'23.25"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = '  +3.02%  '
$ws.Range("D49").Value = "'6.80"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").Value = "'2.21"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = '  -2.66%  '
$ws.Range("E51").Value = '  +0.09%  '
